$p = $ppt.ActivePresentation
$newDate = "2/20/2026"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        }
        if (-not $isDatePlaceholder -and $sh.Name -like "Date Placeholder*") {
            $isDatePlaceholder = $true
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -ne $newDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Update the Date Placeholder on the Slide Master.
$m = $p.SlideMaster
Update-DatePlaceholder $m.Shapes

# Update the Date Placeholder on every slide layout (CustomLayout) tied to the master.
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $cl = $m.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes
}
